$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ORG_PROV_Phone_Scenario")

# Update the cell value from 150 to 1 (cell D2)
$ws.Range("D2").Value = "1"

# Update the active selection on the sheet
$ws.Range("D21").Select()
